$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Thermocouple DIN Rail Terminal Blocks ---
$ws.Range("A10").Value = "Thermocouple DIN Rail Terminal Blocks, Narrow 10.7 mm Width"
$ws.Range("B10").Value = "https://www.omega.com/en-us/temperature-measurement/temperature-connectors-panels-and-block-assemblies/terminal-blocks-and-lugs/drtb-2/p/DRTB-T-2"

# --- Row 11: Thermocouple DIN Terminal Block Connector ---
$ws.Range("A11").Value = "Thermocouple DIN Terminal Block Connector 2 Position Feed Through Beige 12-26 AWG"
$ws.Range("B11").Value = "https://www.digikey.ca/en/products/detail/weidm%C3%BCller/1024100000/497593"

# --- Row 12: DIN Rail Thermocouple Input Signal Conditioners ---
$ws.Range("A12").Value = "DIN Rail Thermocouple Input Signal Conditioners | Low Profile"
$ws.Range("B12").Value = "https://www.omega.com/en-us/data-acquisition/signal-conditioners/din-rail-signal-conditioners/p/DRSL-TC-Srs-Sig-Cond"

# --- Row 13: STATUS SEM1605/TC temp transmitter ---
$ws.Range("A13").Value = "STATUS SEM1605/TC, TEMP TRANSMITTER, THERMOCOUPLE, DIN RAIL"
$ws.Range("B13").Value = "https://www.newark.com/status/sem1605-tc/temp-transmitter-thermocouple/dp/13AC9411?MER=TARG-MER-PDP-RECO-STM71168"

# --- Row 14: SENECA WK109TC0 signal conditioner ---
$ws.Range("A14").Value = "SENECA WK109TC0 SIGNAL CONDITIONER, FOR THEROCOUPLES"
$ws.Range("B14").Value = "https://www.newark.com/seneca/wk109tc0/signal-conditioner-for-therocouples/dp/24M9179"

# --- Row 15: Portenta Machine Control ---
$ws.Range("A15").Value = "Portenta Machine Control"
$ws.Range("B15").Value = "https://store-usa.arduino.cc/products/arduino-portenta-machine-control?selectedStore=us"

# Wrap text + taller rows for the two long entries (A13 / A14)
$ws.Range("A13").WrapText = $true
$ws.Range("A14").WrapText = $true
$ws.Rows("13").RowHeight = 30
$ws.Rows("14").RowHeight = 30

# Hyperlinks for all the newly added B cells (B10:B15)
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.omega.com/en-us/temperature-measurement/temperature-connectors-panels-and-block-assemblies/terminal-blocks-and-lugs/drtb-2/p/DRTB-T-2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.digikey.ca/en/products/detail/weidm%C3%BCller/1024100000/497593") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.omega.com/en-us/data-acquisition/signal-conditioners/din-rail-signal-conditioners/p/DRSL-TC-Srs-Sig-Cond") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.newark.com/status/sem1605-tc/temp-transmitter-thermocouple/dp/13AC9411?MER=TARG-MER-PDP-RECO-STM71168") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.newark.com/seneca/wk109tc0/signal-conditioner-for-therocouples/dp/24M9179") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B15"), "https://store-usa.arduino.cc/products/arduino-portenta-machine-control?selectedStore=us") | Out-Null

# Restore the Hyperlink cell style (Hyperlinks.Add applies its own transient
# formatting) so these cells match the existing hyperlink-styled column B cells
$ws.Range("B10:B15").Style = "Hyperlink"

# Move the selection to where the user ended up after typing the last entry
$ws.Range("B16").Select() | Out-Null
